# Weekly refresh of the "Espárragos" (Vega Monumental Concepción) price
# report: the per-record observations in rows 2-21 are rotated so that
# each row now shows the figures that used to belong to a different row
# (dates, volumes, prices, unit of sale and origin), as published in the
# newer weekly export of the data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2;  D=44875; H='Sin especificar'; I='Primera'; J=300;  K=1500; L=1600; M=1550; N='$/kilo';    O='Provincia de Linares'; P=1550},
    @{Row=3;  D=44519; H='Sin especificar'; I='Primera'; J=250;  K=1200; L=1300; M=1240; N='$/kilo';    O='Provincia de Linares'; P=1240},
    @{Row=4;  D=44511; H='Sin especificar'; I='Primera'; J=600;  K=1300; L=1400; M=1350; N='$/kilo';    O='Provincia de Linares'; P=1350},
    @{Row=5;  D=44489; H='Sin especificar'; I='Primera'; J=600;  K=1400; L=1500; M=1450; N='$/kilo';    O='Provincia de Linares'; P=1450},
    @{Row=6;  D=44876; H='Sin especificar'; I='Primera'; J=350;  K=1500; L=1600; M=1557; N='$/kilo';    O='Provincia de Linares'; P=1557},
    @{Row=7;  D=44510; H='Sin especificar'; I='Primera'; J=600;  K=1300; L=1400; M=1350; N='$/kilo';    O='Provincia de Linares'; P=1350},
    @{Row=8;  D=45202; H='Verde';           I='Primera'; J=300;  K=1600; L=1600; M=1600; N='$/kilo';    O='Provincia de Linares'; P=1600},
    @{Row=9;  D=44860; H='Sin especificar'; I='Primera'; J=1100; K=1500; L=1700; M=1609; N='$/kilo';    O='Provincia de Linares'; P=1609},
    @{Row=10; D=44881; H='Sin especificar'; I='Primera'; J=200;  K=2600; L=2700; M=2650; N='$/kilo';    O='Provincia de Linares'; P=2650},
    @{Row=11; D=44881; H='Sin especificar'; I='Segunda'; J=100;  K=2400; L=2400; M=2400; N='$/kilo';    O='Provincia de Linares'; P=2400},
    @{Row=12; D=44545; H='Sin especificar'; I='Primera'; J=550;  K=1700; L=1800; M=1755; N='$/kilo';    O='Provincia de Linares'; P=1755},
    @{Row=13; D=44496; H='Sin especificar'; I='Primera'; J=550;  K=1500; L=2000; M=1773; N='$/paquete'; O='Provincia de Linares'; P=1773},
    @{Row=14; D=44524; H='Sin especificar'; I='Primera'; J=200;  K=1500; L=1600; M=1550; N='$/kilo';    O='Provincia de Talca';   P=1550},
    @{Row=15; D=44526; H='Sin especificar'; I='Primera'; J=100;  K=1500; L=1600; M=1550; N='$/kilo';    O='Provincia de Linares'; P=1550},
    @{Row=16; D=44839; H='Sin especificar'; I='Primera'; J=500;  K=1700; L=1800; M=1760; N='$/kilo';    O='Provincia de Linares'; P=1760},
    @{Row=17; D=44477; H='Sin especificar'; I='Primera'; J=500;  K=1400; L=1500; M=1460; N='$/kilo';    O='Provincia de Linares'; P=1460},
    @{Row=18; D=45203; H='Sin especificar'; I='Primera'; J=400;  K=1400; L=1500; M=1450; N='$/kilo';    O='Provincia de Linares'; P=1450},
    @{Row=19; D=44468; H='Verde';           I='Primera'; J=500;  K=1800; L=2000; M=1920; N='$/kilo';    O='Provincia de Linares'; P=1920},
    @{Row=20; D=44868; H='Sin especificar'; I='Primera'; J=1000; K=1200; L=1300; M=1250; N='$/kilo';    O='Región del Maule';     P=1250},
    @{Row=21; D=44868; H='Sin especificar'; I='Segunda'; J=200;  K=1000; L=1000; M=1000; N='$/kilo';    O='Región del Maule';     P=1000}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D    # D: Fecha
    $ws.Cells.Item($r, 8).Value  = $u.H    # H: Variedad
    $ws.Cells.Item($r, 9).Value  = $u.I    # I: Calidad
    $ws.Cells.Item($r, 10).Value = $u.J    # J: Volumen
    $ws.Cells.Item($r, 11).Value = $u.K    # K: Precio mínimo
    $ws.Cells.Item($r, 12).Value = $u.L    # L: Precio máximo
    $ws.Cells.Item($r, 13).Value = $u.M    # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $u.N    # N: Unidad de comercialización
    $ws.Cells.Item($r, 15).Value = $u.O    # O: Origen
    $ws.Cells.Item($r, 16).Value = $u.P    # P: Precio $/Kg
}
